# Update vm_pu.xlsx results for the 380 kV case: bus 2 (column B) target voltage
# setpoint changes from 1.05 p.u. to 1.02 p.u., which also changes the computed
# per-unit voltage magnitudes for all other buses/time steps in rows 2:25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F for rows 2:25 (one row of data per time step)
$bf = New-Object "object[,]" 24,5
# row 2 (time step 0)
$bf[0,0] = 1.02
$bf[0,1] = 1.033704477849932
$bf[0,2] = 1.040454095084078
$bf[0,3] = 1.037347739083005
$bf[0,4] = 1.049369600783932
# row 3 (time step 1)
$bf[1,0] = 1.02
$bf[1,1] = 1.034744829903986
$bf[1,2] = 1.041226758718673
$bf[1,3] = 1.038336371803252
$bf[1,4] = 1.050283604881784
# row 4 (time step 2)
$bf[2,0] = 1.02
$bf[2,1] = 1.035418416309605
$bf[2,2] = 1.041726893472858
$bf[2,3] = 1.038976844590581
$bf[2,4] = 1.050875513135181
# row 5 (time step 3)
$bf[3,0] = 1.02
$bf[3,1] = 1.035701689674384
$bf[3,2] = 1.04193718974466
$bf[3,3] = 1.039246280863105
$bf[3,4] = 1.051124466812794
# row 6 (time step 4)
$bf[4,0] = 1.02
$bf[4,1] = 1.035749258236136
$bf[4,2] = 1.041972501708105
$bf[4,3] = 1.039291531021201
$bf[4,4] = 1.051166273971079
# row 7 (time step 5)
$bf[5,0] = 1.02
$bf[5,1] = 1.035422201038863
$bf[5,2] = 1.041729703307519
$bf[5,3] = 1.038980444098951
$bf[5,4] = 1.05087883921405
# row 8 (time step 6)
$bf[6,0] = 1.02
$bf[6,1] = 1.034055984823216
$bf[6,2] = 1.040715184629612
$bf[6,3] = 1.037681694056214
$bf[6,4] = 1.049678391406539
# row 9 (time step 7)
$bf[7,0] = 1.02
$bf[7,1] = 1.031651682801451
$bf[7,2] = 1.03892881731193
$bf[7,3] = 1.035398994045857
$bf[7,4] = 1.047566819644705
# row 10 (time step 8)
$bf[8,0] = 1.02
$bf[8,1] = 1.030050944233162
$bf[8,2] = 1.037738868172419
$bf[8,3] = 1.033881178768442
$bf[8,4] = 1.0461616992223
# row 11 (time step 9)
$bf[9,0] = 1.02
$bf[9,1] = 1.029358314123711
$bf[9,2] = 1.03722384729763
$bf[9,3] = 1.033224901501573
$bf[9,4] = 1.045553893807326
# row 12 (time step 10)
$bf[10,0] = 1.02
$bf[10,1] = 1.029101115706081
$bf[10,2] = 1.037032581829761
$bf[10,3] = 1.032981273754749
$bf[10,4] = 1.045328222021445
# row 13 (time step 11)
$bf[11,0] = 1.02
$bf[11,1] = 1.029156282233327
$bf[11,2] = 1.037073607268303
$bf[11,3] = 1.033033526259336
$bf[11,4] = 1.04537662510787
# row 14 (time step 12)
$bf[12,0] = 1.02
$bf[12,1] = 1.029337052486791
$bf[12,2] = 1.037208036490287
$bf[12,3] = 1.03320476023837
$bf[12,4] = 1.045535237780034
# row 15 (time step 13)
$bf[13,0] = 1.02
$bf[13,1] = 1.029448440979431
$bf[13,2] = 1.037290867597112
$bf[13,3] = 1.033310282091706
$bf[13,4] = 1.045632976792896
# row 16 (time step 14)
$bf[14,0] = 1.02
$bf[14,1] = 1.030096922359469
$bf[14,2] = 1.037773053440218
$bf[14,3] = 1.033924753728515
$bf[14,4] = 1.046202050427316
# row 17 (time step 15)
$bf[15,0] = 1.02
$bf[15,1] = 1.030503831588475
$bf[15,2] = 1.03807557943216
$bf[15,3] = 1.034310449458
$bf[15,4] = 1.046559182184841
# row 18 (time step 16)
$bf[16,0] = 1.02
$bf[16,1] = 1.03074122311853
$bf[16,2] = 1.038252060333225
$bf[16,3] = 1.034535510559544
$bf[16,4] = 1.04676755095332
# row 19 (time step 17)
$bf[17,0] = 1.02
$bf[17,1] = 1.030822175678891
$bf[17,2] = 1.038312239566665
$bf[17,3] = 1.034612266047957
$bf[17,4] = 1.046838609443679
# row 20 (time step 18)
$bf[18,0] = 1.02
$bf[18,1] = 1.030460169086556
$bf[18,2] = 1.038043118911149
$bf[18,3] = 1.034269058507526
$bf[18,4] = 1.046520859129256
# row 21 (time step 19)
$bf[19,0] = 1.02
$bf[19,1] = 1.029283818074887
$bf[19,2] = 1.037168449427863
$bf[19,3] = 1.033154332151144
$bf[19,4] = 1.045488527693495
# row 22 (time step 20)
$bf[20,0] = 1.02
$bf[20,1] = 1.028544634892505
$bf[20,2] = 1.036618719998315
$bf[20,3] = 1.032454286242269
$bf[20,4] = 1.044840005454605
# row 23 (time step 21)
$bf[21,0] = 1.02
$bf[21,1] = 1.02893644858312
$bf[21,2] = 1.036910121759466
$bf[21,3] = 1.032825315248924
$bf[21,4] = 1.045183747352199
# row 24 (time step 22)
$bf[22,0] = 1.02
$bf[22,1] = 1.03047989814251
$bf[22,2] = 1.038057786353542
$bf[22,3] = 1.034287761014254
$bf[22,4] = 1.046538175483571
# row 25 (time step 23)
$bf[23,0] = 1.02
$bf[23,1] = 1.03227287826072
$bf[23,2] = 1.039390470741744
$bf[23,3] = 1.035988427493246
$bf[23,4] = 1.048112258648245

$ws.Range("B2:F25").Value = $bf

# Columns I:N for rows 2:25 (column H stays blank, column G stays 1)
$inArr = New-Object "object[,]" 24,6
# row 2 (time step 0)
$inArr[0,0] = 1.035357102243254
$inArr[0,1] = 1.038827742878673
$inArr[0,2] = 1.043236398255807
$inArr[0,3] = 1.04013888146677
$inArr[0,4] = 1.052126849879535
$inArr[0,5] = 1.040302997428768
# row 3 (time step 1)
$inArr[1,0] = 1.035534362075435
$inArr[1,1] = 1.039510459682119
$inArr[1,2] = 1.043819587881565
$inArr[1,3] = 1.040936829677618
$inArr[1,4] = 1.052852823457393
$inArr[1,5] = 1.040986683768384
# row 4 (time step 2)
$inArr[2,0] = 1.035647802142055
$inArr[2,1] = 1.039952061307971
$inArr[2,2] = 1.044196461196129
$inArr[2,3] = 1.04145330439039
$inArr[2,4] = 1.0533224322373
$inArr[2,5] = 1.041428912519197
# row 5 (time step 3)
$inArr[3,0] = 1.035695190564808
$inArr[3,1] = 1.040137671312515
$inArr[3,2] = 1.044354780871201
$inArr[3,3] = 1.041670465350216
$inArr[3,4] = 1.053519820292201
$inArr[3,5] = 1.041614786111238
# row 6 (time step 4)
$inArr[4,0] = 1.035703129581348
$inArr[4,1] = 1.040168833726378
$inArr[4,2] = 1.044381356508918
$inArr[4,3] = 1.041706929664016
$inArr[4,4] = 1.053552960510373
$inArr[4,5] = 1.041645992779302
# row 7 (time step 5)
$inArr[5,0] = 1.035648436534769
$inArr[5,1] = 1.039954541593163
$inArr[5,2] = 1.044198577135401
$inArr[5,3] = 1.041456205969084
$inArr[5,4] = 1.053325069886604
$inArr[5,5] = 1.041431396326678
# row 8 (time step 6)
$inArr[6,0] = 1.035417268456694
$inArr[6,1] = 1.039058503368613
$inArr[6,2] = 1.04343359082056
$inArr[6,3] = 1.040408520587768
$inArr[6,4] = 1.052372225497585
$inArr[6,5] = 1.040534085625072
# row 9 (time step 7)
$inArr[7,0] = 1.035000296673484
$inArr[7,1] = 1.037478355992436
$inArr[7,2] = 1.042081870418607
$inArr[7,3] = 1.038563533629493
$inArr[7,4] = 1.050692120972176
$inArr[7,5] = 1.038951694258376
# row 10 (time step 8)
$inArr[8,0] = 1.034715867936504
$inArr[8,1] = 1.036424141642261
$inArr[8,2] = 1.041178265458479
$inArr[8,3] = 1.037334364761797
$inArr[8,4] = 1.049571379839915
$inArr[8,5] = 1.037895982802904
# row 11 (time step 9)
$inArr[9,0] = 1.034591182574201
$inArr[9,1] = 1.035967476051817
$inArr[9,2] = 1.04078641976899
$inArr[9,3] = 1.036802323388899
$inArr[9,4] = 1.049085937119783
$inArr[9,5] = 1.037438668694942
# row 12 (time step 10)
$inArr[10,0] = 1.034544639978474
$inArr[10,1] = 1.035797822966018
$inArr[10,2] = 1.040640784488387
$inArr[10,3] = 1.036604729629845
$inArr[10,4] = 1.04890559961558
$inArr[10,5] = 1.0372687746823
# row 13 (time step 11)
$inArr[11,0] = 1.03455463387601
$inArr[11,1] = 1.035834215350055
$inArr[11,2] = 1.040672027647914
$inArr[11,3] = 1.03664711279163
$inArr[11,4] = 1.048944283627113
$inArr[11,5] = 1.037305218747697
# row 14 (time step 12)
$inArr[12,0] = 1.03458734001625
$inArr[12,1] = 1.035953453026032
$inArr[12,2] = 1.040774383269135
$inArr[12,3] = 1.036785989610171
$inArr[12,4] = 1.049071030819424
$inArr[12,5] = 1.037424625754851
# row 15 (time step 13)
$inArr[13,0] = 1.034607461027341
$inArr[13,1] = 1.036026915709501
$inArr[13,2] = 1.040837436522308
$inArr[13,3] = 1.036871560205199
$inArr[13,4] = 1.049149120993159
$inArr[13,5] = 1.037498192763757
# row 16 (time step 14)
$inArr[14,0] = 1.034724110774606
$inArr[14,1] = 1.03645444520023
$inArr[14,2] = 1.041204258824255
$inArr[14,3] = 1.037369678787417
$inArr[14,4] = 1.049603593883906
$inArr[14,5] = 1.037926329395401
# row 17 (time step 15)
$inArr[15,0] = 1.034796873807541
$inArr[15,1] = 1.036722574268114
$inArr[15,2] = 1.041434202238134
$inArr[15,3] = 1.037682188485174
$inArr[15,4] = 1.049888631899171
$inArr[15,5] = 1.038194839237316
# row 18 (time step 16)
$inArr[16,0] = 1.034839168009071
$inArr[16,1] = 1.036878951649815
$inArr[16,2] = 1.041568268456369
$inArr[16,3] = 1.037864489063584
$inArr[16,4] = 1.050054874802932
$inArr[16,5] = 1.038351438692836
# row 19 (time step 17)
$inArr[17,0] = 1.034853564253263
$inArr[17,1] = 1.036932269268092
$inArr[17,2] = 1.041613972074625
$inArr[17,3] = 1.03792665204972
$inArr[17,4] = 1.0501115567792
$inArr[17,5] = 1.038404832028247
# row 20 (time step 18)
$inArr[18,0] = 1.034789082250447
$inArr[18,1] = 1.036693808391014
$inArr[18,2] = 1.041409537261955
$inArr[18,3] = 1.037648657187335
$inArr[18,4] = 1.049858051558493
$inArr[18,5] = 1.038166032509371
# row 21 (time step 19)
$inArr[19,0] = 1.034577715189166
$inArr[19,1] = 1.035918341248179
$inArr[19,2] = 1.040744244468459
$inArr[19,3] = 1.036745093013538
$inArr[19,4] = 1.049033707544574
$inArr[19,5] = 1.037389464114246
# row 22 (time step 20)
$inArr[20,0] = 1.034443496149948
$inArr[20,1] = 1.035430617412739
$inArr[20,2] = 1.0403254492234
$inArr[20,3] = 1.036177160740331
$inArr[20,4] = 1.048515280142107
$inArr[20,5] = 1.036901047655016
# row 23 (time step 21)
$inArr[21,0] = 1.034514773605094
$inArr[21,1] = 1.035689183759354
$inArr[21,2] = 1.040547507626772
$inArr[21,3] = 1.03647821568003
$inArr[21,4] = 1.048790120356216
$inArr[21,5] = 1.037159981195502
# row 24 (time step 22)
$inArr[22,0] = 1.034792603374685
$inArr[22,1] = 1.036706806507541
$inArr[22,2] = 1.041420682475684
$inArr[22,3] = 1.037663808478893
$inArr[22,4] = 1.049871869544337
$inArr[22,5] = 1.038179049084715
# row 25 (time step 23)
$inArr[23,0] = 1.035109231997862
$inArr[23,1] = 1.037887002365135
$inArr[23,2] = 1.042431758784876
$inArr[23,3] = 1.039040364720527
$inArr[23,4] = 1.051126590172152
$inArr[23,5] = 1.039360920955797

$ws.Range("I2:N25").Value = $inArr

Write-Host "Updated vm_pu results for 380 kV case (rows 2:25)."
